{"js": "// Replace every paragraph containing the old Slovak \"Perseus\" observation\n// dates with the new \"Orion\" observation dates, regardless of how the\n// paragraph's text is split across runs.\nconst oldText =\n  \"V roku 2018 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Perseus: 30. okt\u00f3bra - 8. novembra a 29. novembra - 8. decembra\";\nconst newText =\n  \"V roku Orion: 16.-25. janu\u00e1ra, 14.-23. febru\u00e1ra, 14.-24. marca\";\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(oldText) !== -1) {\n    targets.push(paragraphs.items[i]);\n  }\n}\n\nfor (const paragraph of targets) {\n  // Clear all runs/content (and their formatting) from the paragraph,\n  // then insert a brand-new, plain run with the replacement text.\n  paragraph.clear();\n  paragraph.insertText(newText, Word.InsertLocation.start);\n}\n\nawait context.sync();\n", "ps1": "# Replace every paragraph containing the old Slovak \"Perseus\" observation\n# dates with the new \"Orion\" observation dates, regardless of how the\n# paragraph's text is split across runs.\n$d = $word.ActiveDocument\n\n$oldText = \"V roku 2018 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie Perseus: 30. okt\u00f3bra - 8. novembra a 29. novembra - 8. decembra\"\n$newText = \"V roku Orion: 16.-25. janu\u00e1ra, 14.-23. febru\u00e1ra, 14.-24. marca\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $pRng = $p.Range\n    if ($pRng.Text -like \"*$oldText*\") {\n        # Exclude the trailing paragraph mark from the range so we don't\n        # merge with the next paragraph, then wipe all runs/formatting in\n        # the paragraph and insert a brand-new, plain run with the\n        # replacement text.\n        $pRng.End = $pRng.End - 1\n        $pRng.Delete()\n        $pRng.InsertAfter($newText)\n    }\n}\n"}
